# Commit: feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" between "总计" and "2022-Q1", carrying
#    the fund-holding snapshot for that quarter.
# 2. Record the new quarter in the "总计" roll-up sheet (rename its existing
#    row from "2022-Q1" to "2022-Q4" and append a fresh "2022-Q1" row below
#    it, matching the numbers that were already there).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q4" sheet, placed right after "总计" (Add() inserts before
#    the currently-active sheet, i.e. "2022-Q1", which lands it exactly
#    in the middle). Do all sheet lookups-by-name AFTER this structural
#    change so references bind to the post-insert sheet order.
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add()
$wsQ4.Name = "2022-Q4"
# Match the default outline settings the other sheets in this workbook
# already carry (summary row below / summary column on the right).
$wsQ4.Outline.SummaryRow = 1
$wsQ4.Outline.SummaryColumn = -4152

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Match the bordered/bold/centered header style already used on the
# "总计" sheet's header row.
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ4.Range("A2").Value = 0
# Match the "总计" sheet's A2 style (same bordered/bold/centered look).
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B2 and C2..G2 hold plain numeric-looking text in the source data (not
# real numbers). Stamp a text format just long enough to write them so
# Excel doesn't auto-coerce them to numbers, then drop back to the
# workbook's default "Normal" style so no stray formatting is left
# behind on the cell.
$wsQ4.Range("B2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "486002"
$wsQ4.Range("B2").Style = "Normal"

$wsQ4.Range("C2").NumberFormat = "@"
$wsQ4.Range("C2").Value = "工银全球精选股票（QDII）"
$wsQ4.Range("C2").Style = "Normal"

$wsQ4.Range("D2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "3.92"
$wsQ4.Range("D2").Style = "Normal"

$wsQ4.Range("E2").NumberFormat = "@"
$wsQ4.Range("E2").Value = "94.38"
$wsQ4.Range("E2").Style = "Normal"

$wsQ4.Range("F2").NumberFormat = "@"
$wsQ4.Range("F2").Value = "1.73"
$wsQ4.Range("F2").Style = "Normal"

$wsQ4.Range("G2").NumberFormat = "@"
$wsQ4.Range("G2").Value = "0.0678"
$wsQ4.Range("G2").Style = "Normal"

$wsQ4.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 2. "总计" roll-up sheet: the existing data row now describes 2022-Q4,
#    and a new row is appended for 2022-Q1 with the values the sheet
#    used to show.
# ---------------------------------------------------------------------
$wsTotal.Range("B2").Value = "2022-Q4"

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q1"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.07000000000000001

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the originally-active "2022-Q1" sheet selected, same as before
# the edit (adding a sheet otherwise leaves the new one active).
$wsQ1.Activate()

Write-Output "2022-Q4 sheet inserted; 总计 updated"
